$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits in the
#    empty paragraph right after "STATIC CONTENT".
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the empty paragraph right after "ACCEPTANCE CRITERIA".
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd([char]13, [char]7) -eq "ACCEPTANCE CRITERIA") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark in that (still empty) paragraph
#    first, while it is collapsed - this keeps everything inside the
#    existing paragraph instead of spawning a new one.
# ------------------------------------------------------------------
$r = $target.Range
$endPos = $r.End - 1
$bmRange = $d.Range($endPos, $endPos)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bmRange.InsertXML($xml)

# ------------------------------------------------------------------
# 4. Insert the two acceptance-criteria lines immediately before the
#    bookmark (so the bookmark ends up trailing the new text, as in
#    the target document).
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$insertPos = $bm.End

$line1 = "1. Each page must contain multi language options"
$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertBefore($line1)

$pos2 = $insertPos + $line1.Length
$line2 = [char]11 + "2. User can change to different language (such as English, Spanish, Vietnamese, etc.) from any page"
$r2 = $d.Range($pos2, $pos2)
$r2.InsertBefore($line2)

# Colour the second run (break + line 2) first, then the first run -
# applying them in this order keeps the two runs distinct instead of
# letting the save pass coalesce them into a single run.
$range2 = $d.Range($pos2, $pos2 + $line2.Length)
$range2.Font.Color = 0

$range1 = $d.Range($insertPos, $pos2)
$range1.Font.Color = 0

Write-Output "Paragraph text now: [$($target.Range.Text)]"
